$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 6 "Regular US Data" - fill in the input values + formulas
# ---------------------------------------------------------------------------
$ws.Cells.Item(6, 2).Value = 8
$ws.Cells.Item(6, 3).Value = 11
$ws.Cells.Item(6, 4).Value = 27
$ws.Cells.Item(6, 5).Value = 337136173
$ws.Cells.Item(6, 6).Value = 5
$ws.Cells.Item(6, 7).Formula = '=($B$2/B6+$B$2/D6-$B$2/C6)*F6'
$ws.Cells.Item(6, 9).Formula = '=IF(H6>E6,"2ncre","ifc")'

# ---------------------------------------------------------------------------
# Row 7 "Population Increase - High Birth Rate"
# ---------------------------------------------------------------------------
$ws.Cells.Item(7, 2).Value = 10
$ws.Cells.Item(7, 3).Value = 11
$ws.Cells.Item(7, 4).Value = 27
$ws.Cells.Item(7, 5).Value = 337136173
$ws.Cells.Item(7, 6).Value = 5
$ws.Cells.Item(7, 7).Formula = '=($B$2/B7+$B$2/D7-$B$2/C7)*F7'
$ws.Cells.Item(7, 9).Formula = '=IF(H7>E7,"2ncre","ifc")'

# ---------------------------------------------------------------------------
# Row 8 "Population Increase - High Migration"
# ---------------------------------------------------------------------------
$ws.Cells.Item(8, 2).Value = 8
$ws.Cells.Item(8, 3).Value = 11
$ws.Cells.Item(8, 4).Value = 30
$ws.Cells.Item(8, 5).Value = 337136173
$ws.Cells.Item(8, 6).Value = 5
$ws.Cells.Item(8, 7).Formula = '=($B$2/B8+$B$2/D8-$B$2/C8)*F8'
$ws.Cells.Item(8, 9).Formula = '=IF(H8>E8,"2ncre","ifc")'

# ---------------------------------------------------------------------------
# Row 9 "Population Decrease - High Death Rate"
# ---------------------------------------------------------------------------
$ws.Cells.Item(9, 2).Value = 8
$ws.Cells.Item(9, 3).Value = 5
$ws.Cells.Item(9, 4).Value = 27
$ws.Cells.Item(9, 5).Value = 337136173
$ws.Cells.Item(9, 6).Value = 5
$ws.Cells.Item(9, 7).Formula = '=($B$2/B9+$B$2/D9-$B$2/C9)*F9'
$ws.Cells.Item(9, 9).Formula = '=IF(H9>E9,"2ncre","ifc")'

# ---------------------------------------------------------------------------
# Row 10 "Population Low Birth Rate and Low Migration"
# ---------------------------------------------------------------------------
$ws.Cells.Item(10, 2).Value = 5
$ws.Cells.Item(10, 3).Value = 11
$ws.Cells.Item(10, 4).Value = 15
$ws.Cells.Item(10, 5).Value = 337136173
$ws.Cells.Item(10, 6).Value = 5
$ws.Cells.Item(10, 7).Formula = '=($B$2/B10+$B$2/D10-$B$2/C10)*F10'
$ws.Cells.Item(10, 9).Formula = '=IF(H10>E10,"2ncre","ifc")'

# ---------------------------------------------------------------------------
# Row heights / column width
# ---------------------------------------------------------------------------
$ws.Rows.Item(5).RowHeight = 60
$ws.Rows.Item(6).RowHeight = 43.5
$ws.Columns.Item(7).ColumnWidth = 10.5

# ---------------------------------------------------------------------------
# Selection moves from B6:I10 to the single cell G9
# ---------------------------------------------------------------------------
[void]$ws.Range("G9").Select()
